# Applies the "Updated symbol list" edit (cryptos.xlsx, Fri Jan 13 00:55:20 UTC 2023)
# Rewrites Price/Volume(1h)/Data/Hora for each coin row, and for rows whose
# ranking shifted, the Coin name + Link as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D="286.91"; E="0.76%"; F="13-1-2023"; G="0"},
    @{Row=3; D="29.27"; E="3.57%"; F="13-1-2023"; G="0"},
    @{Row=4; D="5.103"; E="3.58%"; F="13-1-2023"; G="0"},
    @{Row=5; D="0.06624"; E="1.46%"; F="13-1-2023"; G="0"},
    @{Row=6; D="7.360"; E="2.82%"; F="13-1-2023"; G="0"},
    @{Row=7; B="GateToken"; C="https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; D="3.411"; E="1.79%"; F="13-1-2023"; G="0"},
    @{Row=8; B="FTXToken"; C="https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; D="1.366"; E="5.30%"; F="13-1-2023"; G="0"},
    @{Row=9; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="0.9270"; E="0.14%"; F="13-1-2023"; G="0"},
    @{Row=10; B="WazirX"; C="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; D="0.1571"; E="0.95%"; F="13-1-2023"; G="0"},
    @{Row=11; B="LiechtensteinCryptoassetsExchange"; C="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; D="0.06536"; E="4.85%"; F="13-1-2023"; G="0"},
    @{Row=12; B="MandalaExchangeToken"; C="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D="0.07619"; E="-0.06%"; F="13-1-2023"; G="0"},
    @{Row=13; B="BitrueCoin"; C="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D="0.02905"; E="-0.73%"; F="13-1-2023"; G="0"},
    @{Row=14; B="BitMartToken"; C="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D="0.08981"; E="-0.15%"; F="13-1-2023"; G="0"},
    @{Row=15; B="BitForexToken"; C="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D="0.001592"; E="-0.57%"; F="13-1-2023"; G="0"},
    @{Row=16; B="CoinExToken"; C="https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; D="0.04496"; E="2.30%"; F="13-1-2023"; G="0"},
    @{Row=17; B="One"; C="https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; D="0.0006455"; E="-0.40%"; F="13-1-2023"; G="0"},
    @{Row=18; B="TigerCash"; C="https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D="0.006287"; E="3.28%"; F="13-1-2023"; G="0"},
    @{Row=19; B="LEO"; C="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D="3.445"; E="-2.06%"; F="13-1-2023"; G="0"},
    @{Row=20; D="2.236"; E="-0.12%"; F="13-1-2023"; G="0"},
    @{Row=21; D="0.3213"; E="2.04%"; F="13-1-2023"; G="0"},
    @{Row=22; E="-4.15%"; F="13-1-2023"; G="0"},
    @{Row=23; D="4.080"; E="2.43%"; F="13-1-2023"; G="0"},
    @{Row=24; D="0.1570"; E="0.25%"; F="13-1-2023"; G="0"},
    @{Row=25; E="1.43%"; F="13-1-2023"; G="0"},
    @{Row=26; D="0.004115"; E="-4.99%"; F="13-1-2023"; G="0"},
    @{Row=27; D="0.0001251"; E="6.79%"; F="13-1-2023"; G="0"},
    @{Row=28; D="0.0001619"; E="-1.58%"; F="13-1-2023"; G="0"},
    @{Row=29; F="13-1-2023"; G="0"},
    @{Row=30; F="13-1-2023"; G="0"},
    @{Row=31; F="13-1-2023"; G="0"},
    @{Row=32; F="13-1-2023"; G="0"},
    @{Row=33; F="13-1-2023"; G="0"},
    @{Row=34; F="13-1-2023"; G="0"},
    @{Row=35; F="13-1-2023"; G="0"},
    @{Row=36; F="13-1-2023"; G="0"},
    @{Row=37; F="13-1-2023"; G="0"},
    @{Row=38; F="13-1-2023"; G="0"},
    @{Row=39; F="13-1-2023"; G="0"},
    @{Row=40; D="0.04183"; E="1.05%"; F="13-1-2023"; G="0"},
    @{Row=41; D="0.006722"; E="4.35%"; F="13-1-2023"; G="0"},
    @{Row=42; D="0.1241"; E="-11.20%"; F="13-1-2023"; G="0"},
    @{Row=43; D="0.001981"; E="-6.74%"; F="13-1-2023"; G="0"},
    @{Row=44; D="0.01246"; E="8.84%"; F="13-1-2023"; G="0"},
    @{Row=45; D="0.00005622"; E="1.65%"; F="13-1-2023"; G="0"},
    @{Row=46; B="CoinbaseStockToken"; C="https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"; D="0.01308"; E="-29.10%"; F="13-1-2023"; G="0"},
    @{Row=47; B="BOLO"; C="https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"; D="1.968"; E="26.05%"; F="13-1-2023"; G="0"},
    @{Row=48; F="13-1-2023"; G="0"},
    @{Row=49; F="13-1-2023"; G="0"},
    @{Row=50; F="13-1-2023"; G="0"},
    @{Row=51; F="13-1-2023"; G="0"}
)

$colMap = @{B=2; C=3; D=4; E=5; F=6; G=7}

foreach ($item in $data) {
    $r = $item.Row
    foreach ($col in @("B","C","D","E","F","G")) {
        if ($item.ContainsKey($col)) {
            $cell = $ws.Cells.Item($r, $colMap[$col])
            $cell.NumberFormat = "@"
            $cell.Value = $item[$col]
        }
    }
}
